$d = $word.ActiveDocument

# Generic helper: locate `old` with Find (no replace), then assign the new
# text directly onto a fresh Range built from the match's Start/End. This
# mirrors a literal "select found text, type new text" flow and keeps the
# run's existing formatting intact.
function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    $fresh = $d.Range($rng.Start, $rng.End)
    $fresh.Text = $new
    return $true
}

# Some runs sit immediately after a closing </w:hyperlink> with no explicit
# run formatting of their own (they just render in the surrounding default
# style). Replacing text starting exactly at that boundary position causes
# the new text to inherit the hyperlink's (blue/underlined) formatting
# instead. Avoid this by leaving the first character of the match untouched
# (it keeps anchoring the run's real formatting) and only replacing the
# remainder of the matched text.
function Replace-TextAfterHyperlink($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    $fresh = $d.Range($rng.Start + 1, $rng.End)
    $fresh.Text = $new.Substring(1)
    return $true
}

# Language switcher line (top of doc) + standalone language label below it.
# "English" occurs twice in the document (once in the hyperlink, once on its
# own further down); each Replace-Text call only touches the next remaining
# match, so call it once per occurrence.
Replace-Text "English" "Inglês"
Replace-TextAfterHyperlink " / Portuguese / French / Thai / Vietnamese / Spanish" " / Português / Francês / Tailandês / Vietnamita / Espanhol"
Replace-Text "English" "Inglês"

# Brief table
Replace-Text "Brief" "Resumo"
Replace-Text "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io" "Um correio eletrónico enviado a parceiros no país alvo cujos documentos não passaram no nosso processo de verificação. Será enviado através de customer.io"
Replace-Text "Target audience" "Público-alvo"
Replace-Text "Invited partners who submitted wrong/incomplete documents" "Parceiros convidados que submeteram documentos incorretos/incompletos"

# Subject line
Replace-Text "Subject line" "Linha de assunto"
Replace-Text "[EVENT NAME]" "[NOME DO EVENTO]"
Replace-Text " — document verification failed " " - A verificação dos documentos falhou "

# Heading
Replace-Text "Uh oh! Your documents couldn’t be verified" "Uh oh! Os seus documentos não puderam ser verificados"

# Greeting
Replace-Text "Hi " "Olá "
Replace-Text "[PARTNER NAME]" "[NOME DO PARCEIRO]"

# Body
Replace-Text "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "Lamentamos informar que os seus documentos não passaram no nosso processo de verificação, uma vez que encontrámos os seguintes problemas: "

# Bulleted list
Replace-Text "A copy of your vaccination certificate" "Uma cópia do seu certificado de vacinação"
Replace-Text ": Document is unclear" ": O documento não é claro"
Replace-Text "[Document 2]" "[Documento 2]"

# Resubmission deadline
Replace-Text "Please resubmit the documents above by " "Por favor, reenvie os documentos acima até "
Replace-Text "DD Mmm YYYY" "DD Mmm AAAA"
Replace-Text " so we can proceed with the necessary arrangements." " para que possamos proceder às devidas diligências."

# Contact via live chat / WhatsApp, and contact country manager.
# " or " occurs twice (after "live chat" and after "[EMAIL ADDRESS]"); the
# first one sits right after the "live chat" hyperlink, so it needs the
# boundary-safe replacement helper.
Replace-Text "If you have any questions, please contact us via " "Para mais informações, contacte-nos através do "
Replace-TextAfterHyperlink " or " " ou "
Replace-Text "If you have any questions, please contact your country manager, " "Para mais questões, pode também contactar o seus gestor de parcerias "
Replace-Text ", at " ", em "
Replace-Text " or " " ou "
